# Information disclosure.docx edit
# Commit: "Ubacen celokupan STRIDE u prezentaciju." (Inserted the whole STRIDE into the presentation)
#
# The only substantive textual change in this particular file is in the
# opening paragraph: the trailing sentence explaining why information
# disclosure matters is removed, leaving just the closing period after
# "... vredne informacije o sistemu."
#
#   "... vredne informacije o sistemu. Zbog toga je neophodno uvek
#    razmotriti koje se informacije otkrivaju i da li mogu biti
#    zloupotrebljene od strane zlonamernog korisnika."
# becomes
#   "... vredne informacije o sistemu."

$d = $word.ActiveDocument

# Word drops a collapsed "_GoBack" bookmark at the point of the last edit
# whenever the document is saved after a change. Mark that spot now --
# right after "... vredne informacije o sistemu." -- before the trailing
# sentence is removed, so it lands exactly where editing will finish.
$goBack = $d.Range(122, 122)
$d.Bookmarks.Add("_GoBack", $goBack)

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = " Zbog toga je neophodno uvek razmotriti koje se informacije otkrivaju i da li mogu biti zloupotrebljene od strane zlonamernog korisnika."
$find.Replacement.ClearFormatting()
$find.Replacement.Text = ""
$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)
